{"js": "// Update the date and the 25 two-digit-by-two-digit multiplication\n// expressions to the new values from the commit (\"output generated at\n// c986bee\"). Each old value is unique within the document, so a simple\n// search+replace per pair is safe and keeps the original run formatting\n// (font/size) intact because insertText(..., Replace) only swaps the\n// text of the matched range.\n\nconst replacements = [\n  [\"2024-10-16 Wednesday\", \"2024-10-17 Thursday\"],\n  [\"84\u00d737=\", \"34\u00d734=\"],\n  [\"33\u00d797=\", \"94\u00d797=\"],\n  [\"17\u00d794=\", \"68\u00d772=\"],\n  [\"51\u00d796=\", \"82\u00d779=\"],\n  [\"61\u00d744=\", \"57\u00d745=\"],\n  [\"76\u00d737=\", \"11\u00d743=\"],\n  [\"49\u00d716=\", \"63\u00d782=\"],\n  [\"67\u00d723=\", \"15\u00d724=\"],\n  [\"87\u00d748=\", \"77\u00d746=\"],\n  [\"57\u00d721=\", \"40\u00d732=\"],\n  [\"41\u00d744=\", \"22\u00d724=\"],\n  [\"88\u00d721=\", \"35\u00d716=\"],\n  [\"69\u00d764=\", \"78\u00d790=\"],\n  [\"41\u00d750=\", \"99\u00d760=\"],\n  [\"77\u00d785=\", \"36\u00d768=\"],\n  [\"30\u00d786=\", \"91\u00d792=\"],\n  [\"80\u00d773=\", \"37\u00d733=\"],\n  [\"18\u00d783=\", \"28\u00d795=\"],\n  [\"17\u00d762=\", \"13\u00d794=\"],\n  [\"55\u00d787=\", \"45\u00d750=\"],\n  [\"36\u00d746=\", \"20\u00d784=\"],\n  [\"84\u00d769=\", \"83\u00d775=\"],\n  [\"91\u00d727=\", \"66\u00d752=\"],\n  [\"86\u00d763=\", \"15\u00d791=\"],\n  [\"32\u00d768=\", \"42\u00d774=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date and the 25 two-digit-by-two-digit multiplication\n# expressions to the new values from the commit (\"output generated at\n# c986bee\"). Each old value is unique within the document, so a plain\n# Find/Replace per pair is safe; Find.Execute only swaps the text inside\n# the matched range, so the run's original formatting (font/size) is\n# left untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    [PSCustomObject]@{ Old = \"2024-10-16 Wednesday\"; New = \"2024-10-17 Thursday\" },\n    [PSCustomObject]@{ Old = \"84\u00d737=\"; New = \"34\u00d734=\" },\n    [PSCustomObject]@{ Old = \"33\u00d797=\"; New = \"94\u00d797=\" },\n    [PSCustomObject]@{ Old = \"17\u00d794=\"; New = \"68\u00d772=\" },\n    [PSCustomObject]@{ Old = \"51\u00d796=\"; New = \"82\u00d779=\" },\n    [PSCustomObject]@{ Old = \"61\u00d744=\"; New = \"57\u00d745=\" },\n    [PSCustomObject]@{ Old = \"76\u00d737=\"; New = \"11\u00d743=\" },\n    [PSCustomObject]@{ Old = \"49\u00d716=\"; New = \"63\u00d782=\" },\n    [PSCustomObject]@{ Old = \"67\u00d723=\"; New = \"15\u00d724=\" },\n    [PSCustomObject]@{ Old = \"87\u00d748=\"; New = \"77\u00d746=\" },\n    [PSCustomObject]@{ Old = \"57\u00d721=\"; New = \"40\u00d732=\" },\n    [PSCustomObject]@{ Old = \"41\u00d744=\"; New = \"22\u00d724=\" },\n    [PSCustomObject]@{ Old = \"88\u00d721=\"; New = \"35\u00d716=\" },\n    [PSCustomObject]@{ Old = \"69\u00d764=\"; New = \"78\u00d790=\" },\n    [PSCustomObject]@{ Old = \"41\u00d750=\"; New = \"99\u00d760=\" },\n    [PSCustomObject]@{ Old = \"77\u00d785=\"; New = \"36\u00d768=\" },\n    [PSCustomObject]@{ Old = \"30\u00d786=\"; New = \"91\u00d792=\" },\n    [PSCustomObject]@{ Old = \"80\u00d773=\"; New = \"37\u00d733=\" },\n    [PSCustomObject]@{ Old = \"18\u00d783=\"; New = \"28\u00d795=\" },\n    [PSCustomObject]@{ Old = \"17\u00d762=\"; New = \"13\u00d794=\" },\n    [PSCustomObject]@{ Old = \"55\u00d787=\"; New = \"45\u00d750=\" },\n    [PSCustomObject]@{ Old = \"36\u00d746=\"; New = \"20\u00d784=\" },\n    [PSCustomObject]@{ Old = \"84\u00d769=\"; New = \"83\u00d775=\" },\n    [PSCustomObject]@{ Old = \"91\u00d727=\"; New = \"66\u00d752=\" },\n    [PSCustomObject]@{ Old = \"86\u00d763=\"; New = \"15\u00d791=\" },\n    [PSCustomObject]@{ Old = \"32\u00d768=\"; New = \"42\u00d774=\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
